# Daily attendance processing - 2026-02-21 11:04:09 UTC
# Reorders the comma-separated values in column G ("Recorded By") for a
# number of rows in the "Session Analysis Results" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$updates = @{
    3  = "2022/2023, 2025/2026"
    22 = "2024/2025, 2025/2026"
    23 = "2022/2023, 2025/2026, 2023/2024"
    24 = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    27 = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    28 = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    31 = "2022/2023, 2025/2026"
    50 = "2024/2025, 2025/2026"
    51 = "2022/2023, 2025/2026, 2023/2024"
    52 = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    55 = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    56 = "2025/2026, neveen.nashaat@med.asu.edu.eg"
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
